$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.20253495330645
$ws.Range("C2").Value = 11.9335494769093
$ws.Range("E2").Value = 13.68832009907383
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 31.38951654343337
$ws.Range("H2").Value = 15.15826671466092
$ws.Range("K2").Value = 8.028348091289596
$ws.Range("L2").Value = 9.913115098108733
$ws.Range("M2").Value = 13.48900533334152
$ws.Range("O2").Value = 23.38049157789171
$ws.Range("B3").Value = 10.9369276873992
$ws.Range("C3").Value = 11.96382260459514
$ws.Range("E3").Value = 13.72950882435556
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 31.58621115992892
$ws.Range("H3").Value = 15.21777881630419
$ws.Range("K3").Value = 7.775666420954377
$ws.Range("L3").Value = 9.919515703729276
$ws.Range("M3").Value = 13.44434674713366
$ws.Range("O3").Value = 23.49471621886505
$ws.Range("B4").Value = 10.77176071625128
$ws.Range("C4").Value = 11.98366052179299
$ws.Range("E4").Value = 13.75692596179812
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 31.71566761722013
$ws.Range("H4").Value = 15.25646425507106
$ws.Range("K4").Value = 7.615028896840967
$ws.Range("L4").Value = 9.924735486835157
$ws.Range("M4").Value = 13.41844692531842
$ws.Range("O4").Value = 23.56924385787204
$ws.Range("B5").Value = 10.70402323515211
$ws.Range("C5").Value = 11.99205967908321
$ws.Range("E5").Value = 13.76863377133275
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 31.77060199380627
$ws.Range("H5").Value = 15.27276908272756
$ws.Range("K5").Value = 7.548251799935884
$ws.Range("L5").Value = 9.927187572936726
$ws.Range("M5").Value = 13.40828182978879
$ws.Range("O5").Value = 23.6007198270683
$ws.Range("B6").Value = 10.69275227611341
$ws.Range("C6").Value = 11.99347340130679
$ws.Range("E6").Value = 13.77061016969118
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 31.77985534221793
$ws.Range("H6").Value = 15.27550913956632
$ws.Range("K6").Value = 7.537086036944554
$ws.Range("L6").Value = 9.927614385969882
$ws.Range("M6").Value = 13.40661763046264
$ws.Range("O6").Value = 23.60601316770303
$ws.Range("B7").Value = 10.77084880215648
$ws.Range("C7").Value = 11.98377251922572
$ws.Range("E7").Value = 13.75708169049378
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 31.7163996608098
$ws.Range("H7").Value = 15.25668195926603
$ws.Range("K7").Value = 7.614133556247607
$ws.Range("L7").Value = 9.924767239864725
$ws.Range("M7").Value = 13.41830825020242
$ws.Range("O7").Value = 23.56966387740312
$ws.Range("B8").Value = 11.11144048216358
$ws.Range("C8").Value = 11.94372861375906
$ws.Range("E8").Value = 13.70208063407878
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 31.45553173111302
$ws.Range("H8").Value = 15.17834186899038
$ws.Range("K8").Value = 7.942399678535267
$ws.Range("L8").Value = 9.915054742431883
$ws.Range("M8").Value = 13.47329565632551
$ws.Range("O8").Value = 23.41896479249364
$ws.Range("B9").Value = 11.75894481751098
$ws.Range("C9").Value = 11.8750903180181
$ws.Range("E9").Value = 13.61109081590604
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 31.01307086599556
$ws.Range("H9").Value = 15.04169485843252
$ws.Range("K9").Value = 8.540044650353314
$ws.Range("L9").Value = 9.906211923414237
$ws.Range("M9").Value = 13.59287637586853
$ws.Range("O9").Value = 23.15827389377937
$ws.Range("B10").Value = 12.21740382048131
$ws.Range("C10").Value = 11.83064668671011
$ws.Range("E10").Value = 13.55450852302315
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 30.73038094280168
$ws.Range("H10").Value = 14.9515941922274
$ws.Range("K10").Value = 8.948072620581401
$ws.Range("L10").Value = 9.905891769086841
$ws.Range("M10").Value = 13.68745475212214
$ws.Range("O10").Value = 22.98793699694395
$ws.Range("B11").Value = 12.42130471564595
$ws.Range("C11").Value = 11.81171868955428
$ws.Range("E11").Value = 13.53099408713453
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 30.61104042299435
$ws.Range("H11").Value = 14.91282859747461
$ws.Range("K11").Value = 9.12643158587275
$ws.Range("L11").Value = 9.907076672278494
$ws.Range("M11").Value = 13.73184065211609
$ws.Range("O11").Value = 22.91504050900808
$ws.Range("B12").Value = 12.49777848460114
$ws.Range("C12").Value = 11.80473591083228
$ws.Range("E12").Value = 13.52240950188967
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 30.56718564701397
$ws.Range("H12").Value = 14.89846772066513
$ws.Range("K12").Value = 9.192891414218387
$ws.Range("L12").Value = 9.907715652810541
$ws.Range("M12").Value = 13.7488353444979
$ws.Range("O12").Value = 22.88809627367735
$ws.Range("B13").Value = 12.48134241284153
$ws.Range("C13").Value = 11.80623156577937
$ws.Range("E13").Value = 13.5242441243122
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 30.57657100523181
$ws.Range("H13").Value = 14.90154642342528
$ws.Range("K13").Value = 9.17862672260912
$ws.Range("L13").Value = 9.907569591496051
$ws.Range("M13").Value = 13.74516707659131
$ws.Range("O13").Value = 22.8938698395958
$ws.Range("B14").Value = 12.42761145077871
$ws.Range("C14").Value = 11.8111405110548
$ws.Range("E14").Value = 13.53028141934939
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 30.60740563717265
$ws.Range("H14").Value = 14.9116407337154
$ws.Range("K14").Value = 9.131921136939038
$ws.Range("L14").Value = 9.90712543429524
$ws.Range("M14").Value = 13.73323512659069
$ws.Range("O14").Value = 22.91281056156218
$ws.Range("B15").Value = 12.394601459469
$ws.Range("C15").Value = 11.81417143597302
$ws.Range("E15").Value = 13.5340210823844
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 30.6264670167199
$ws.Range("H15").Value = 14.91786528985822
$ws.Range("K15").Value = 9.103170821128179
$ws.Range("L15").Value = 9.906878123331696
$ws.Range("M15").Value = 13.72595050443886
$ws.Range("O15").Value = 22.92449826257644
$ws.Range("B16").Value = 12.20397866741067
$ws.Range("C16").Value = 11.83190957470238
$ws.Range("E16").Value = 13.55609000953078
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 30.73836688519119
$ws.Range("H16").Value = 14.95417226359314
$ws.Range("K16").Value = 8.936267137370521
$ws.Range("L16").Value = 9.905841019032803
$ws.Range("M16").Value = 13.68458065825347
$ws.Range("O16").Value = 22.99279332394516
$ws.Range("B17").Value = 12.08579424309042
$ws.Range("C17").Value = 11.84312122176476
$ws.Range("E17").Value = 13.57019838862526
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 30.8093887424535
$ws.Range("H17").Value = 14.97701395943404
$ws.Range("K17").Value = 8.831990760543539
$ws.Range("L17").Value = 9.905544852236886
$ws.Range("M17").Value = 13.6595439714322
$ws.Range("O17").Value = 23.03586585191917
$ws.Range("B18").Value = 12.01738341108643
$ws.Range("C18").Value = 11.84969128673668
$ws.Range("E18").Value = 13.57852260494011
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 30.85110954992011
$ws.Range("H18").Value = 14.99036104905541
$ws.Range("K18").Value = 8.771332815872174
$ws.Range("L18").Value = 9.905499809382246
$ws.Range("M18").Value = 13.64527226366374
$ws.Range("O18").Value = 23.06107212083004
$ws.Range("B19").Value = 11.99414830848001
$ws.Range("C19").Value = 11.85193667075782
$ws.Range("E19").Value = 13.58137701912612
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 30.86538490399938
$ws.Range("H19").Value = 14.99491608641479
$ws.Range("K19").Value = 8.750679347322988
$ws.Range("L19").Value = 9.905506110464337
$ws.Range("M19").Value = 13.64046249177229
$ws.Range("O19").Value = 23.06968074449176
$ws.Range("B20").Value = 12.09842064590777
$ws.Range("C20").Value = 11.84191516046772
$ws.Range("E20").Value = 13.56867485136837
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 30.80173816656736
$ws.Range("H20").Value = 14.97456078047211
$ws.Range("K20").Value = 8.843161919196058
$ws.Range("L20").Value = 9.905563418503954
$ws.Range("M20").Value = 13.66219591131949
$ws.Range("O20").Value = 23.03123598857151
$ws.Range("B21").Value = 12.44341412368214
$ws.Range("C21").Value = 11.80969362302993
$ws.Range("E21").Value = 13.52849944085316
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 30.59831242857193
$ws.Range("H21").Value = 14.90866714275641
$ws.Range("K21").Value = 9.145669310599478
$ws.Range("L21").Value = 9.907250738515931
$ws.Range("M21").Value = 13.73673483582593
$ws.Range("O21").Value = 22.90722929800719
$ws.Range("B22").Value = 12.6645509343237
$ws.Range("C22").Value = 11.78971211672626
$ws.Range("E22").Value = 13.50410657201031
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 30.47315715311079
$ws.Range("H22").Value = 14.86745983399063
$ws.Range("K22").Value = 9.337061064636915
$ws.Range("L22").Value = 9.90946207763734
$ws.Range("M22").Value = 13.78653426937489
$ws.Range("O22").Value = 22.83003124318533
$ws.Range("B23").Value = 12.54694445959632
$ws.Range("C23").Value = 11.80027826303169
$ws.Range("E23").Value = 13.5169549942773
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 30.53923966246969
$ws.Range("H23").Value = 14.88928315597429
$ws.Range("K23").Value = 9.235500717486191
$ws.Range("L23").Value = 9.908180774913841
$ws.Range("M23").Value = 13.75985926002963
$ws.Range("O23").Value = 22.87088125248163
$ws.Range("B24").Value = 12.0927136923061
$ws.Range("C24").Value = 11.84246003343204
$ws.Range("E24").Value = 13.56936297870674
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 30.80519422196328
$ws.Range("H24").Value = 14.97566919262472
$ws.Range("K24").Value = 8.838113640961765
$ws.Range("L24").Value = 9.905554634550228
$ws.Range("M24").Value = 13.66099658791667
$ws.Range("O24").Value = 23.0333277689404
$ws.Range("B25").Value = 11.58648013556261
$ws.Range("C25").Value = 11.89260471120577
$ws.Range("E25").Value = 13.63390171712309
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 31.12534419097798
$ws.Range("H25").Value = 15.07684979523311
$ws.Range("K25").Value = 8.383627459912448
$ws.Range("L25").Value = 9.907516014282418
$ws.Range("M25").Value = 13.55931294853462
$ws.Range("O25").Value = 23.22507331456106
